# Auto-generated Excel COM-interop script to apply scheduled-runner price updates
# to the Siren_Profits-style Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 499.85715
$ws.Range("I39").Value = 333
$ws.Range("K39").Value = 999
$ws.Range("M39").Value = -703

$ws.Range("H40").Value = 5399.2
$ws.Range("I40").Value = 5000
$ws.Range("J40").Value = 5499
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 5499
$ws.Range("M40").Value = -4825
$ws.Range("N40").Value = -5849

$ws.Range("H94").Value = 50120252
$ws.Range("I94").Value = 71428890
$ws.Range("K94").Value = 71428890
$ws.Range("M94").Value = -71428439

$ws.Range("H100").Value = 99929.664
$ws.Range("I100").Value = 111916.11
$ws.Range("K100").Value = 111916.11
$ws.Range("M100").Value = -111375.11

$ws.Range("H137").Value = 8244.5
$ws.Range("I137").Value = 13486.353
$ws.Range("K137").Value = 40459.05899999999
$ws.Range("M137").Value = -37909.05899999999

$ws.Range("H141").Value = 5754.033
$ws.Range("I141").Value = 5565.522
$ws.Range("K141").Value = 16696.566
$ws.Range("M141").Value = -11516.566

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6527.5693
$ws.Range("I32").Value = 6359.5483
$ws.Range("K32").Value = 6359.5483
$ws.Range("M32").Value = -6072.5483

$ws.Range("H61").Value = 9043.963
$ws.Range("J61").Value = 5244.5
$ws.Range("L61").Value = 5244.5
$ws.Range("N61").Value = -5668.5

$ws.Range("H74").Value = 5830.5557
$ws.Range("I74").Value = 6339.2856
$ws.Range("J74").Value = 4050
$ws.Range("K74").Value = 6339.2856
$ws.Range("L74").Value = 4050
$ws.Range("M74").Value = -5465.2856
$ws.Range("N74").Value = -5798

$ws.Range("H77").Value = 5830.5557
$ws.Range("I77").Value = 6339.2856
$ws.Range("J77").Value = 4050
$ws.Range("K77").Value = 31696.428
$ws.Range("L77").Value = 20250
$ws.Range("M77").Value = -27328.428
$ws.Range("N77").Value = -28986

$ws.Range("H132").Value = 4363.5415
$ws.Range("I132").Value = 4232.263
$ws.Range("J132").Value = 4862.4
$ws.Range("K132").Value = 12696.789
$ws.Range("L132").Value = 14587.2
$ws.Range("M132").Value = -10166.789
$ws.Range("N132").Value = -19647.2

$ws.Range("H136").Value = 9043.963
$ws.Range("J136").Value = 5244.5
$ws.Range("L136").Value = 15733.5
$ws.Range("N136").Value = -20833.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H94").Value = 7806.8535
$ws.Range("I94").Value = 9632.467000000001
$ws.Range("J94").Value = 2827.9092
$ws.Range("K94").Value = 9632.467000000001
$ws.Range("L94").Value = 2827.9092
$ws.Range("M94").Value = -9181.467000000001
$ws.Range("N94").Value = -3729.9092

$ws.Range("H105").Value = 56571.367
$ws.Range("I105").Value = 79442.92
$ws.Range("J105").Value = 7016.3335
$ws.Range("K105").Value = 79442.92
$ws.Range("L105").Value = 7016.3335
$ws.Range("M105").Value = -77695.92
$ws.Range("N105").Value = -10510.3335

$ws.Range("H134").Value = 19368.285
$ws.Range("I134").Value = 19368.285
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 58104.855
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = -55569.855
$ws.Range("M134").ClearContents()

$ws.Range("H140").Value = 94743
$ws.Range("J140").Value = 94743
$ws.Range("L140").Value = 94743
$ws.Range("N140").Value = -105103

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1266.5
$ws.Range("I22").Value = 866.3333
$ws.Range("K22").Value = 866.3333
$ws.Range("M22").Value = -516.3333

$ws.Range("H31").Value = 10187.947
$ws.Range("I31").Value = 12507.363
$ws.Range("J31").Value = 6998.75
$ws.Range("K31").Value = 12507.363
$ws.Range("L31").Value = 6998.75
$ws.Range("M31").Value = -12212.363
$ws.Range("N31").Value = -7588.75

$ws.Range("H34").Value = 10187.947
$ws.Range("I34").Value = 12507.363
$ws.Range("J34").Value = 6998.75
$ws.Range("K34").Value = 12507.363
$ws.Range("L34").Value = 6998.75
$ws.Range("M34").Value = -12305.363
$ws.Range("N34").Value = -7402.75

$ws.Range("H99").Value = 10561828
$ws.Range("I99").Value = 29031278
$ws.Range("K99").Value = 29031278
$ws.Range("M99").Value = -29029780

$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -44920

$ws.Range("H126").Value = 10561828
$ws.Range("I126").Value = 29031278
$ws.Range("K126").Value = 87093834
$ws.Range("M126").Value = -87091364

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 261.57144
$ws.Range("I33").Value = 500
$ws.Range("J33").Value = 221.83333
$ws.Range("K33").Value = 3000
$ws.Range("L33").Value = 1330.99998
$ws.Range("M33").Value = -2717
$ws.Range("N33").Value = -1896.99998

$ws.Range("H56").Value = 6207.5
$ws.Range("I56").Value = 6207.5
$ws.Range("K56").Value = 6207.5
$ws.Range("M56").Value = -5677.5

$ws.Range("H122").Value = 4880.2705
$ws.Range("J122").Value = 5270.909
$ws.Range("L122").Value = 47438.181
$ws.Range("N122").Value = -52338.181

$ws.Range("H140").Value = 14656.846
$ws.Range("I140").Value = 15482.417
$ws.Range("J140").Value = 4750
$ws.Range("K140").Value = 46447.251
$ws.Range("L140").Value = 14250
$ws.Range("M140").Value = -41267.251
$ws.Range("N140").Value = -24610

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 92825.75
$ws.Range("J134").Value = 92825.75
$ws.Range("L134").Value = 278477.25
$ws.Range("N134").Value = -283547.25

$ws.Range("H141").Value = 82846.71000000001
$ws.Range("J141").Value = 82846.71000000001
$ws.Range("L141").Value = 82846.71000000001
$ws.Range("N141").Value = -93206.71000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 19177.756
$ws.Range("I40").Value = 21775.738
$ws.Range("K40").Value = 21775.738
$ws.Range("M40").Value = -21639.738

$ws.Range("H132").Value = 1244363.9
$ws.Range("I132").Value = 1863896.5
$ws.Range("K132").Value = 5591689.5
$ws.Range("M132").Value = -5589159.5

$ws.Range("H140").Value = 130485.8
$ws.Range("J140").Value = 130485.8
$ws.Range("L140").Value = 130485.8
$ws.Range("N140").Value = -140845.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4116.5684
$ws.Range("I122").Value = 2344.111
$ws.Range("K122").Value = 7032.333
$ws.Range("M122").Value = -4582.333

$ws.Range("H132").Value = 9687.927
$ws.Range("I132").Value = 11550.518
$ws.Range("J132").Value = 5186.6665
$ws.Range("K132").Value = 34651.554
$ws.Range("L132").Value = 15559.9995
$ws.Range("M132").Value = -32121.554
$ws.Range("N132").Value = -20619.9995

Write-Host "Applied scheduled-runner price/profit updates to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR."
